$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (System Name) becomes a plain numeric sequence 1..12,
# replacing the old text labels (Nokia AMS, Viavi Fusion, ...).
# Column B (Date) values are left untouched.
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# Update the viewport selection to match the new view state.
$ws.Range("D16").Select()
